$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 154, shifting existing rows 154:221 down to 155:222.
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row 154 with the new record.
$ws.Range("A154").Value2 = 10
$ws.Range("B154").Value2 = "Vega Modelo de Temuco"
$ws.Range("C154").Value2 = "La Araucanía"
$ws.Range("D154").Value2 = 45141
$ws.Range("E154").Value2 = 9
$ws.Range("F154").Value2 = 100112031
$ws.Range("G154").Value2 = "Poroto verde"
$ws.Range("H154").Value2 = "Sin especificar"
$ws.Range("I154").Value2 = "Primera"
$ws.Range("J154").Value2 = 80
$ws.Range("K154").Value2 = 32000
$ws.Range("L154").Value2 = 32000
$ws.Range("M154").Value2 = 32000
$ws.Range("N154").Value2 = "`$/malla 25 kilos"
$ws.Range("O154").Value2 = "Región de Arica y Parinacota"
$ws.Range("P154").Value2 = 1280
$ws.Range("Q154").Value2 = 25
$ws.Range("R154").Value2 = "Hortaliza"
